$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, borders, center/top alignment) from an
# existing header cell (H1) onto the two new header cells.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Data values for columns I (I0) and J (IF), rows 2-20
$values = @(
    @(9, 9),
    @(2, 4),
    @(8, 8),
    @(4, 5),
    @(10, 11),
    @(8, 9),
    @(9, 9),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 7),
    @(5, 8),
    @(4, 8),
    @(5, 9),
    @(6, 7),
    @(1, 3),
    @(7, 8),
    @(5, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
